$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 7 ("input_efficiency"), pushing
# everything down by two rows. Formatting of the rows below carries down
# automatically; the two freshly inserted rows start out blank.
$ws.Rows("7:8").Insert()

# Fill in the shared text first, in the same order the original authoring
# tool introduced it (input, output, configuration_fxe), then the rest.
$ws.Cells.Item(7, 3).Value2 = "input"
$ws.Cells.Item(8, 3).Value2 = "output"
$ws.Cells.Item(7, 4).Value2 = "configuration_fxe"
$ws.Cells.Item(8, 4).Value2 = "configuration_fxe"

# New row 7: "input" / "configuration_fxe" parameter for the hydro flow.
$ws.Cells.Item(7, 1).Value2 = "CHE"
$ws.Cells.Item(7, 2).Value2 = "conv_elec_hydroror"
$ws.Cells.Item(7, 6).Value2 = "hydro"
$ws.Cells.Item(7, 7).Value2 = 1

# New row 8: "output" / "configuration_fxe" parameter for the elecsupply flow.
$ws.Cells.Item(8, 1).Value2 = "CHE"
$ws.Cells.Item(8, 2).Value2 = "conv_elec_hydroror"
$ws.Cells.Item(8, 6).Value2 = "elecsupply"
$ws.Cells.Item(8, 7).Value2 = 1

# Keep the autofilter / named range in sync with the two extra rows.
$ws.AutoFilterMode = $false
$ws.Range("A5:L852").AutoFilter()

$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$5:`$L`$852"

# Move the active selection the way the authored workbook shows it.
$ws.Range("G9").Select()
